$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 19-57 (columns A-R). Each inner array holds values in column order A..R.
$data = @(
    @(5, "Macroferia Regional de Talca", "Maule", 44434, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 300, 6000, 6000, 6000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 100, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44280, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 7000, 7000, 7000, "`$/caja 60 unidades", "Región del Maule", 117, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44299, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44242, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 300, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44258, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44243, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 300, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44421, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44237, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44273, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44257, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 150, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44295, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44298, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44274, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 150, 7000, 7000, 7000, "`$/caja 60 unidades", "Región del Maule", 117, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44435, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 300, 6000, 6000, 6000, "`$/caja 50 unidades", "Región de Arica y Parinacota", 120, 50, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44435, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 1300, 6000, 7000, 6231, "`$/caja 60 unidades", "Región de Arica y Parinacota", 104, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44431, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 400, 6000, 6000, 6000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 100, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44231, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 250, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44278, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44389, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 300, 12000, 12000, 12000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 200, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44251, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44250, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44305, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44294, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 9000, 9000, 9000, "`$/caja 60 unidades", "Región del Maule", 150, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44417, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 300, 7000, 7000, 7000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 117, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44419, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 7000, 7000, 7000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 117, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44245, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44265, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44277, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44433, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 300, 7000, 7000, 7000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 117, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44309, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 150, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44253, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44272, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44230, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 9000, 9000, 9000, "`$/caja 60 unidades", "Región del Maule", 150, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44232, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 150, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44270, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44244, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44284, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región del Maule", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44418, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 200, 8000, 8000, 8000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 133, 60, "Hortaliza"),
    @(5, "Macroferia Regional de Talca", "Maule", 44432, 7, 100112001, "Berenjena", "Sin especificar", "Primera", 300, 6000, 6000, 6000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 100, 60, "Hortaliza"),
)

$startRow = 19
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
$dateFormat = $ws.Range("D18").NumberFormat

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $col = $cols[$j]
        $ws.Range("$col$r").Value = $rowVals[$j]
    }
    $ws.Range("D$r").NumberFormat = $dateFormat
}

Write-Output "done"